$wb = $excel.ActiveWorkbook

# "registerInfo" is the first sheet (rId1 / sheet1.xml)
$wsRegister = $wb.Worksheets.Item("registerInfo")

# Fill in the "done" column (P) for rows 2-7 with "X", matching the
# pre-existing pattern already present in rows 8-14 of that column.
$wsRegister.Range("P2:P7").Value = "X"

# Copy the cell formatting used by the rest of that data (column G uses the
# style we need) onto the newly filled cells so they pick up style index 5.
$wsRegister.Range("G2").Copy()
$wsRegister.Range("P2:P7").PasteSpecial(-4122)

# Make "registerInfo" the active/selected sheet (it was "testcase1" before),
# with P2 selected as the active cell.
$wsRegister.Activate() | Out-Null
$wsRegister.Range("P2").Select() | Out-Null
